$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 28 and 29 swap coin identity: NEARProtocol <-> FirstDigitalUSD,
# each also getting a new Price/Volume value.
$ws.Range("B28").Value = "FirstDigitalUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.16%  "

$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.27"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +3.06%  "

# Remaining Price (D) / Volume(1h) (E) refresh across the table.
$ws.Range("D2").Value = "63.744.24"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "3.327.30"
$ws.Range("E3").Value = "  +5.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.82"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +2.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.43"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +2.82%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.326.42"
$ws.Range("E8").Value = "  +5.30%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.52"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +4.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.469"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.82"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("D15").Value = "3.864.63"
$ws.Range("E15").Value = "  +5.11%  "
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "3.320.31"
$ws.Range("E17").Value = "  +5.20%  "
$ws.Range("D18").Value = "63.781.49"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.88"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +3.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "480.40"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.737"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +5.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.22"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +6.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.72"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +5.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.94"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +2.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.21"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +3.38%  "
$ws.Range("E31").Value = "  +3.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.07"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +8.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.107"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("E36").Value = "  +4.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.44"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +6.23%  "
$ws.Range("E39").Value = "  +3.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "434.87"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +4.01%  "
$ws.Range("D41").Value = "3.105.92"
$ws.Range("E41").Value = "  +5.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.122"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +9.51%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.35"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.80%  "
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("E46").Value = "  +4.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "37.20"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +15.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.37"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.31"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +2.66%  "
